$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: clear all content (was a DAQLab / low / "Implement a more rigorous way..." entry)
$ws.Range("A3:F3").ClearContents()
$ws.Rows.Item(3).AutoFit()

# Row 4: update issue description text (column D)
$ws.Range("D4").Value = "Handle DAQmx errors when setting up DAQmx tasks and couple them to the task controller state. If an error is encountered while setting up a DAQmx task, then the task controller should be switched to an Unconfigured state so that it cannot be executed from the task tree"

# Row 5: update assigned-to, module and issue description
$ws.Range("B5").Value = "low"
$ws.Range("C5").Value = "DAQLab"
$ws.Range("D5").Value = "When a task controller is dropped as a child of another task controller in the Task Tree, the selection jumps back to the first element in the Task Tree. This is because the task tree is re-assembled each time an a task controller item is dragged and dropped. This gives a slightly annoying user experience when assembing task trees."

# Row 6: new date, module, and issue description (laser scanning / DAQmx interaction)
$ws.Range("A6").Value = 42168
$ws.Range("C6").Value = "Laser scanning and DAQmx interaction"
$ws.Range("D6").Value = "When AI sampling rate is 200 KHz and 3x oversampling is used, so that an actual sampling rate of 600 KHz is used for a 5us pixel dwell time, then the image looks bad. For other values of the oversampling such as 1,2,4,5 the image looks good.  Investigate why other values for the pixel dwell time such as 3.125 us give a bad image. Is the algorithm going wrong somewhere ??"
$ws.Rows.Item(6).RowHeight = 60

# Row 7: replaced with a new issue (data storage / Rawdata default directory)
$ws.Range("A7").Value = 42168
$ws.Range("B7").Value = "low"
$ws.Range("C7").Value = "data storage"
$ws.Range("D7").Value = "If Rawdata default directory is not present, then create it by default."
$ws.Rows.Item(7).AutoFit()

# Row 8: clear all content (was a "high" / DAQLab, task controller / UITC race condition entry)
$ws.Range("A8:F8").ClearContents()
$ws.Rows.Item(8).AutoFit()

# Update the active selection shown when the workbook is opened
$ws.Range("E7").Select()
